$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 393, shifting existing rows 393:490 down to 394:491.
$ws.Rows.Item(393).Insert()

# Populate the newly inserted row 393 with the new weekly data point.
$ws.Range("A393").Value = 9
$ws.Range("B393").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C393").Value = "Metropolitana"
$ws.Range("D393").Value = 44932
$ws.Range("E393").Value = 13
$ws.Range("F393").Value = 100112044
$ws.Range("G393").Value = "Perejil"
$ws.Range("H393").Value = "Sin especificar"
$ws.Range("I393").Value = "Primera"
$ws.Range("J393").Value = 70
$ws.Range("K393").Value = 12000
$ws.Range("L393").Value = 15000
$ws.Range("M393").Value = 13500
$ws.Range("N393").Value = "$/docena de atados"
$ws.Range("O393").Value = "Región Metropolitana"
$ws.Range("P393").Value = 4500
$ws.Range("Q393").Value = 3
$ws.Range("R393").Value = "Hortaliza"

# Match the source style (date format) used by the rest of column D.
$ws.Range("D393").NumberFormat = $ws.Range("D394").NumberFormat
